# Replace the single data row (row 2) on Sheet1 with the freshly
# extracted record. A handful of columns are now blank (F, I, V) while
# a few previously-blank columns are now populated (B, D, S).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dataRange = $ws.Range("A2:AC2")

# Temporarily force a text number format across the row so that
# numeric/date-looking values (dates, zip codes, phone numbers, ids,
# etc.) are written as literal text instead of being auto-converted to
# numbers or dates by Excel's type inference.
$dataRange.NumberFormat = "@"

$ws.Range("A2").Value = "Third-Party Bill"
$ws.Range("B2").Value = "2024-03-05"
$ws.Range("C2").Value = "1039163"
$ws.Range("D2").Value = "16567886"
$ws.Range("E2").Value = "WED31500"
$ws.Range("F2").Value = $null
$ws.Range("G2").Value = "Christian"
$ws.Range("H2").Value = "William"
$ws.Range("I2").Value = $null
$ws.Range("J2").Value = "1964-10-30"
$ws.Range("K2").Value = "Male"
$ws.Range("L2").Value = "11386 68th St N"
$ws.Range("M2").Value = "FL"
$ws.Range("N2").Value = "West Palm Beach"
$ws.Range("O2").Value = "33412"
$ws.Range("P2").Value = "5617236746"
$ws.Range("Q2").Value = "table"
$ws.Range("R2").Value = "Jennifer Marshall, PA-C,"
$ws.Range("S2").Value = "Wellington WED"
$ws.Range("T2").Value = "William Christian"
$ws.Range("U2").Value = "Self"
$ws.Range("V2").Value = $null
$ws.Range("W2").Value = "BCBS of FL Blueselar PPO/EPO/POS/FEP/PPC"
$ws.Range("X2").Value = "99999U6K"
$ws.Range("Y2").Value = "VMAH45391700"
$ws.Range("Z2").Value = "PO BOX 1798"
$ws.Range("AA2").Value = "Jacksonville"
$ws.Range("AB2").Value = "FL"
$ws.Range("AC2").Value = "322310014"

# Restore the default ("Normal") cell style now that the text values
# are locked in, so the cells keep referencing the workbook's original
# (default) style index rather than a newly-forced text style.
$dataRange.Style = "Normal"
